# daily auto push: 2025-10-09 09:29 UTC
# Append the new daily-ranking row (row 85) to the bottom of the sheet.
#
# Column A holds dates formatted/entered as plain text ("2025/10/09"),
# not real Excel date serials, so we seed the new row by copying the
# last existing row (which already carries that literal text plus the
# sheet's default/no-style formatting) straight down, then only touch
# the one cell whose value actually changes. This avoids Excel's
# automatic "looks like a date" reinterpretation that a direct
# `.Value = "2025/10/09"` assignment on a blank General-formatted cell
# would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry A85:D85 down from A84:D84 (date text, weekday text, and the
# ranking value are unchanged for this new entry).
$ws.Range("A84:D84").Copy($ws.Range("A85:D85"))

# Only the hour/time column differs for the new observation.
$ws.Cells.Item(85, 3).Value = 17
